# Weekly update: a new Perejil (Vega Central Mapocho de Santiago) price
# record is inserted as the new row 229, pushing the existing rows
# 229-238 down to 230-239 (dimension grows from A1:R238 to A1:R239).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 229, shifting rows 229:238 down to 230:239.
$ws.Rows.Item(229).Insert()

# Populate the new row 229 with the latest weekly observation.
$ws.Cells.Item(229, 1).Value  = 9
$ws.Cells.Item(229, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(229, 3).Value  = "Metropolitana"
$ws.Cells.Item(229, 4).Value  = 44509
$ws.Cells.Item(229, 5).Value  = 13
$ws.Cells.Item(229, 6).Value  = 100112044
$ws.Cells.Item(229, 7).Value  = "Perejil"
$ws.Cells.Item(229, 8).Value  = "Sin especificar"
$ws.Cells.Item(229, 9).Value  = "Primera"
$ws.Cells.Item(229, 10).Value = 106
$ws.Cells.Item(229, 11).Value = 8000
$ws.Cells.Item(229, 12).Value = 10000
$ws.Cells.Item(229, 13).Value = 9000
$ws.Cells.Item(229, 14).Value = "`$/docena de atados"
$ws.Cells.Item(229, 15).Value = "Región Metropolitana"
$ws.Cells.Item(229, 16).Value = 3000
$ws.Cells.Item(229, 17).Value = 3
$ws.Cells.Item(229, 18).Value = "Hortaliza"
